# Flight.xlsx edit: drop the "Test Type" (Functional/Regression) column from
# Test_Case_List and rename the still-manual test case IDs (TC_001..TC_009)
# to the shorter TC_01..TC_09 form used for cases that need no test data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Case_List")

# Rename TestCase_ID values TC_001..TC_009 -> TC_01..TC_09 (rows 2-10, column B)
for ($i = 2; $i -le 10; $i++) {
    $cell = $ws.Cells.Item($i, 2)
    $cell.Value2 = ($cell.Value2 -replace "TC_00", "TC_0")
}

# Remove the now-unused "Test Type" column (column E); cells shift left
$ws.Range("E1").EntireColumn.Delete()

# Resize the old "Run" column (now D) - no longer a best-fit width
$ws.Columns.Item(4).ColumnWidth = 13

# Reapply the autofilter over the new, narrower range
$ws.AutoFilterMode = $false
$ws.Range("A1:E16").AutoFilter() | Out-Null

# Point the (hidden) _FilterDatabase defined name at the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Test_Case_List!_FilterDatabase") {
        $n.RefersTo = "=Test_Case_List!`$A`$1:`$E`$16"
    }
}

# Move the active selection
$ws.Range("E3").Select() | Out-Null
